$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the GDP observations (rows 12-30) with the latest FRED vintage
$ws.Range("B12").Value = 115056.79399999999
$ws.Range("B13").Value = 122438.054
$ws.Range("B14").Value = 131143.09099999999
$ws.Range("B15").Value = 141858.951
$ws.Range("B16").Value = 152878.644
$ws.Range("B17").Value = 160656.91399999999
$ws.Range("B18").Value = 167416.33300000001
$ws.Range("B19").Value = 168273.109
$ws.Range("B20").Value = 164417.69399999999
$ws.Range("B21").Value = 167899.89199999999
$ws.Range("B22").Value = 173264.40100000001
$ws.Range("B23").Value = 178241.24299999999
$ws.Range("B24").Value = 188473.15
$ws.Range("B25").Value = 197247.96400000001
$ws.Range("B26").Value = 207514.33300000001
$ws.Range("B27").Value = 211862.37400000001
$ws.Range("B28").Value = 222238.592
$ws.Range("B29").Value = 232579.152
$ws.Range("B30").Value = 242950.04199999999

# Append the new observation for 2020-01-01
$ws.Range("A31").Value = 43831
$ws.Range("A31").NumberFormat = "yyyy\-mm\-dd"
$ws.Range("B31").Value = 240411.07199999999
$ws.Range("B31").NumberFormat = "0.000"

# Match the saved view state (whole-column selection over A:B)
$ws.Columns("A:B").Select()
